$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update grouped activity codes in column I to use combined labels
$ws.Range("I35:I39").Value = "DE"
$ws.Range("I46:I50").Value = "HJ"
$ws.Range("I53:I58").Value = "HJ"
$ws.Range("I62:I75").Value = "LMN"
$ws.Range("I81:I89").Value = "RS"

# Rows 82-89 got an explicit font style applied in the source edit
$ws.Range("I82:I89").Font.Name = "Calibri"

# Adjust the visible view/selection like in the authored workbook
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("I82:I89").Select()
